$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 59.7
$ws.Range("N2").Value = 54.83846622768671

# Row 3
$ws.Range("D3").Value = 92351.97
$ws.Range("E3").Value = 62
$ws.Range("F3").Value = 1.65
$ws.Range("K3").Value = 55.7
$ws.Range("N3").Value = 54.83846622768671

# Row 4
$ws.Range("K4").Value = 51.5
$ws.Range("N4").Value = 54.83846622768671

# Row 5
$ws.Range("K5").Value = 49.7
$ws.Range("N5").Value = 54.83846622768671

# Row 6
$ws.Range("K6").Value = 35.9
$ws.Range("N6").Value = 54.83846622768671
